# Updated cryptos list values (prices + 1h volume %) per the latest data pull,
# including the Aptos / InternetComputer(DFINITY) row swap (rows 42-43).
#
# Values that look like plain decimal numbers (single "." separator, e.g. "1.002")
# would otherwise be auto-coerced to numeric cells by Excel, so for those we
# briefly force Text format, write the literal string, then clear the format
# again so the cell keeps the original (default/General) style - matching how
# the source cells were already stored as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "28.220.00"; Text = $false },
    @{ Cell = "E2"; Value = "  +0.06%  "; Text = $false },
    @{ Cell = "D3"; Value = "1.906.64"; Text = $false },
    @{ Cell = "E3"; Value = "  +1.86%  "; Text = $false },
    @{ Cell = "D4"; Value = "1.002"; Text = $true },
    @{ Cell = "E4"; Value = "  -0.06%  "; Text = $false },
    @{ Cell = "D5"; Value = "314.17"; Text = $true },
    @{ Cell = "E5"; Value = "  +0.74%  "; Text = $false },
    @{ Cell = "D6"; Value = "1.002"; Text = $true },
    @{ Cell = "E6"; Value = "  +0.11%  "; Text = $false },
    @{ Cell = "D7"; Value = "0.5088"; Text = $true },
    @{ Cell = "E7"; Value = "  +1.35%  "; Text = $false },
    @{ Cell = "D8"; Value = "0.3932"; Text = $true },
    @{ Cell = "E8"; Value = "  +0.01%  "; Text = $false },
    @{ Cell = "D9"; Value = "0.09640"; Text = $true },
    @{ Cell = "E9"; Value = "  -1.67%  "; Text = $false },
    @{ Cell = "D10"; Value = "1.137"; Text = $true },
    @{ Cell = "E10"; Value = "  +0.02%  "; Text = $false },
    @{ Cell = "D11"; Value = "41.97"; Text = $true },
    @{ Cell = "E11"; Value = "  +1.76%  "; Text = $false },
    @{ Cell = "D12"; Value = "6.420"; Text = $true },
    @{ Cell = "E12"; Value = "  -0.78%  "; Text = $false },
    @{ Cell = "D13"; Value = "20.93"; Text = $true },
    @{ Cell = "E13"; Value = "  -0.39%  "; Text = $false },
    @{ Cell = "D14"; Value = "1.918.97"; Text = $false },
    @{ Cell = "E14"; Value = "  +2.69%  "; Text = $false },
    @{ Cell = "D15"; Value = "7.331"; Text = $true },
    @{ Cell = "E15"; Value = "  -0.93%  "; Text = $false },
    @{ Cell = "E16"; Value = "  -0.03%  "; Text = $false },
    @{ Cell = "D17"; Value = "0.00001122"; Text = $true },
    @{ Cell = "D18"; Value = "92.67"; Text = $true },
    @{ Cell = "E18"; Value = "  -0.86%  "; Text = $false },
    @{ Cell = "D19"; Value = "0.06643"; Text = $true },
    @{ Cell = "E19"; Value = "  +0.23%  "; Text = $false },
    @{ Cell = "D20"; Value = "17.97"; Text = $true },
    @{ Cell = "E21"; Value = "  +0.08%  "; Text = $false },
    @{ Cell = "D22"; Value = "6.237"; Text = $true },
    @{ Cell = "E22"; Value = "  +1.56%  "; Text = $false },
    @{ Cell = "D23"; Value = "28.273.95"; Text = $false },
    @{ Cell = "E23"; Value = "  -0.01%  "; Text = $false },
    @{ Cell = "E24"; Value = "  -0.17%  "; Text = $false },
    @{ Cell = "D25"; Value = "2.309"; Text = $true },
    @{ Cell = "E25"; Value = "  +1.76%  "; Text = $false },
    @{ Cell = "D26"; Value = "2.665"; Text = $true },
    @{ Cell = "E26"; Value = "  +4.02%  "; Text = $false },
    @{ Cell = "D27"; Value = "2.142.52"; Text = $false },
    @{ Cell = "E27"; Value = "  +2.75%  "; Text = $false },
    @{ Cell = "D28"; Value = "20.98"; Text = $true },
    @{ Cell = "E28"; Value = "  -1.62%  "; Text = $false },
    @{ Cell = "D29"; Value = "158.07"; Text = $true },
    @{ Cell = "E29"; Value = "  +0.13%  "; Text = $false },
    @{ Cell = "D30"; Value = "127.08"; Text = $true },
    @{ Cell = "E30"; Value = "  -0.60%  "; Text = $false },
    @{ Cell = "D31"; Value = "1.093"; Text = $true },
    @{ Cell = "E31"; Value = "  +2.67%  "; Text = $false },
    @{ Cell = "E32"; Value = "  +0.03%  "; Text = $false },
    @{ Cell = "D33"; Value = "5.635"; Text = $true },
    @{ Cell = "E33"; Value = "  +0.10%  "; Text = $false },
    @{ Cell = "D34"; Value = "3.624"; Text = $true },
    @{ Cell = "E34"; Value = "  +0.10%  "; Text = $false },
    @{ Cell = "D35"; Value = "9.635"; Text = $true },
    @{ Cell = "E35"; Value = "  +1.34%  "; Text = $false },
    @{ Cell = "D36"; Value = "0.06654"; Text = $true },
    @{ Cell = "E36"; Value = "  -2.04%  "; Text = $false },
    @{ Cell = "D37"; Value = "0.02418"; Text = $true },
    @{ Cell = "E37"; Value = "  +1.22%  "; Text = $false },
    @{ Cell = "D38"; Value = "1.243"; Text = $true },
    @{ Cell = "E38"; Value = "  +1.86%  "; Text = $false },
    @{ Cell = "E39"; Value = "  +0.33%  "; Text = $false },
    @{ Cell = "E40"; Value = "  +11.03%  "; Text = $false },
    @{ Cell = "D41"; Value = "0.6383"; Text = $true },
    @{ Cell = "E41"; Value = "  +1.37%  "; Text = $false },
    @{ Cell = "B42"; Value = "InternetComputer(DFINITY)"; Text = $false },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; Text = $false },
    @{ Cell = "D42"; Value = "5.014"; Text = $true },
    @{ Cell = "E42"; Value = "  -0.06%  "; Text = $false },
    @{ Cell = "B43"; Value = "Aptos"; Text = $false },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; Text = $false },
    @{ Cell = "D43"; Value = "11.46"; Text = $true },
    @{ Cell = "E43"; Value = "  -0.22%  "; Text = $false },
    @{ Cell = "E44"; Value = "  +0.12%  "; Text = $false },
    @{ Cell = "D45"; Value = "13.43"; Text = $true },
    @{ Cell = "E45"; Value = "  -1.45%  "; Text = $false },
    @{ Cell = "D46"; Value = "0.6019"; Text = $true },
    @{ Cell = "E46"; Value = "  +0.02%  "; Text = $false },
    @{ Cell = "D47"; Value = "3.750"; Text = $true },
    @{ Cell = "E47"; Value = "  +2.31%  "; Text = $false },
    @{ Cell = "D48"; Value = "1.282"; Text = $true },
    @{ Cell = "E48"; Value = "  +0.95%  "; Text = $false },
    @{ Cell = "D49"; Value = "2.036"; Text = $true },
    @{ Cell = "E49"; Value = "  +2.31%  "; Text = $false },
    @{ Cell = "D50"; Value = "123.39"; Text = $true },
    @{ Cell = "E50"; Value = "  -1.00%  "; Text = $false },
    @{ Cell = "E51"; Value = "  -0.87%  "; Text = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Text) {
        $range.NumberFormat = "@"
        $range.Value2 = $u.Value
        $range.ClearFormats()
    } else {
        $range.Value2 = $u.Value
    }
}
